# Inserts two new rows (new market-report entries) above the current row 14,
# pushing all existing rows from 14 downward by two rows. Then fills the
# newly inserted row 14/15 with the new "Primera"/"Segunda" price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 14 (existing row 14 -> row 16, etc.)
$ws.Rows("14:15").Insert()

# ---- New row 14 ("Primera") ----
$ws.Cells.Item(14, 1).Value2 = 11
$ws.Cells.Item(14, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(14, 3).Value2 = "Bíobío"
$ws.Cells.Item(14, 4).Value2 = 44490
$ws.Cells.Item(14, 5).Value2 = 8
$ws.Cells.Item(14, 6).Value2 = 100114013
$ws.Cells.Item(14, 7).Value2 = "Zanahoria"
$ws.Cells.Item(14, 8).Value2 = "Sin especificar"
$ws.Cells.Item(14, 9).Value2 = "Primera"
$ws.Cells.Item(14, 10).Value2 = 600
$ws.Cells.Item(14, 11).Value2 = 7000
$ws.Cells.Item(14, 12).Value2 = 7500
$ws.Cells.Item(14, 13).Value2 = 7250
$ws.Cells.Item(14, 14).Value2 = "$/saco 20 kilos"
$ws.Cells.Item(14, 15).Value2 = "Chillán"
$ws.Cells.Item(14, 16).Value2 = 362
$ws.Cells.Item(14, 17).Value2 = 20
$ws.Cells.Item(14, 18).Value2 = "Hortaliza"

# ---- New row 15 ("Segunda") ----
$ws.Cells.Item(15, 1).Value2 = 11
$ws.Cells.Item(15, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(15, 3).Value2 = "Bíobío"
$ws.Cells.Item(15, 4).Value2 = 44490
$ws.Cells.Item(15, 5).Value2 = 8
$ws.Cells.Item(15, 6).Value2 = 100114013
$ws.Cells.Item(15, 7).Value2 = "Zanahoria"
$ws.Cells.Item(15, 8).Value2 = "Sin especificar"
$ws.Cells.Item(15, 9).Value2 = "Segunda"
$ws.Cells.Item(15, 10).Value2 = 300
$ws.Cells.Item(15, 11).Value2 = 6000
$ws.Cells.Item(15, 12).Value2 = 6000
$ws.Cells.Item(15, 13).Value2 = 6000
$ws.Cells.Item(15, 14).Value2 = "$/saco 20 kilos"
$ws.Cells.Item(15, 15).Value2 = "Chillán"
$ws.Cells.Item(15, 16).Value2 = 300
$ws.Cells.Item(15, 17).Value2 = 20
$ws.Cells.Item(15, 18).Value2 = "Hortaliza"

# Give the new date cells the same date style as the rest of column D
$ws.Range("D14:D15").NumberFormat = $ws.Cells.Item(16, 4).NumberFormat
